$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D8").Value = 18095300
$ws.Range("E8").Value = 15720800
$ws.Range("F8").Value = 15440100
$ws.Range("G8").Value = 13714800
$ws.Range("H8").Value = 12901100
$ws.Range("I8").Value = 13083500
$ws.Range("J8").Value = 12246000
$ws.Range("D9").Value = 1835800
$ws.Range("E9").Value = 1728800
$ws.Range("F9").Value = 1777500
$ws.Range("G9").Value = 1643900
$ws.Range("H9").Value = 1484000
$ws.Range("I9").Value = 5681200
$ws.Range("J9").Value = 12536100
$ws.Range("D10").Value = 16259500
$ws.Range("E10").Value = 13992000
$ws.Range("F10").Value = 13662700
$ws.Range("G10").Value = 12071000
$ws.Range("H10").Value = 11417100
$ws.Range("I10").Value = 7402300
$ws.Range("J10").Value = -290000
$ws.Range("D14").Value = 16900
$ws.Range("E14").Value = 26600
$ws.Range("F14").Value = 42700
$ws.Range("G14").Value = 119300
$ws.Range("H14").Value = 88400
$ws.Range("I14").Value = 41300
$ws.Range("D15").Value = 1021000
$ws.Range("E15").Value = 1187600
$ws.Range("F15").Value = 1087600
$ws.Range("G15").Value = 1149000
$ws.Range("H15").Value = 746200
$ws.Range("I15").Value = 729100
$ws.Range("J15").Value = 620900
$ws.Range("D17").Value = 13900500
$ws.Range("E17").Value = 11708700
$ws.Range("F17").Value = 11682700
$ws.Range("G17").Value = 11207700
$ws.Range("H17").Value = 9552900
$ws.Range("I17").Value = 9452300
$ws.Range("J17").Value = 8326500
$ws.Range("D18").Value = 4194800
$ws.Range("E18").Value = 4012100
$ws.Range("F18").Value = 3757500
$ws.Range("G18").Value = 2507100
$ws.Range("H18").Value = 3348200
$ws.Range("I18").Value = 3631200
$ws.Range("J18").Value = 3919500
$ws.Range("D20").Value = -152900
$ws.Range("E20").Value = -4987900
$ws.Range("F20").Value = 10647300
$ws.Range("G20").Value = 14073700
$ws.Range("H20").Value = 2033000
$ws.Range("I20").Value = -24200
$ws.Range("J20").Value = 1422600
$ws.Range("D21").Value = 5064200
$ws.Range("E21").Value = 213400
$ws.Range("F21").Value = 15493800
$ws.Range("G21").Value = 17731300
$ws.Range("H21").Value = 6128400
$ws.Range("I21").Value = 4337000
$ws.Range("J21").Value = 5963900
$ws.Range("D22").Value = 246800
$ws.Range("E22").Value = 136100
$ws.Range("F22").Value = 244100
$ws.Range("G22").Value = 109900
$ws.Range("H22").Value = 77400
$ws.Range("I22").Value = 81800
$ws.Range("J22").Value = 48200
$ws.Range("D23").Value = 3795100
$ws.Range("E23").Value = -1111800
$ws.Range("F23").Value = 14160600
$ws.Range("G23").Value = 16470900
$ws.Range("H23").Value = 5303800
$ws.Range("I23").Value = 3525200
$ws.Range("J23").Value = 5293900
$ws.Range("D24").Value = 796400
$ws.Range("E24").Value = -156500
$ws.Range("F24").Value = 2432400
$ws.Range("G24").Value = 2844500
$ws.Range("H24").Value = 1005900
$ws.Range("I24").Value = 686500
$ws.Range("J24").Value = 1059000
$ws.Range("D26").Value = 2998700
$ws.Range("E26").Value = -955300
$ws.Range("F26").Value = 11728200
$ws.Range("G26").Value = 13626400
$ws.Range("H26").Value = 4297800
$ws.Range("I26").Value = 2838700
$ws.Range("J26").Value = 4234800
$ws.Range("D27").Value = 2997700
$ws.Range("E27").Value = -956300
$ws.Range("F27").Value = 10907200
$ws.Range("G27").Value = 12652500
$ws.Range("H27").Value = 4017800
$ws.Range("I27").Value = 2662700
$ws.Range("J27").Value = 3979700
$ws.Range("D32").Value = 152900
$ws.Range("E32").Value = 4987900
$ws.Range("F32").Value = -10647300
$ws.Range("G32").Value = -14073700
$ws.Range("H32").Value = -2033000
$ws.Range("I32").Value = 24200
$ws.Range("J32").Value = -1422600
$ws.Range("D33").Value = 2997700
$ws.Range("E33").Value = -956300
$ws.Range("F33").Value = 10907200
$ws.Range("G33").Value = 12652500
$ws.Range("H33").Value = 4017800
$ws.Range("I33").Value = 2662700
$ws.Range("J33").Value = 3979700
$ws.Range("D35").Value = 2997700
$ws.Range("E35").Value = -956300
$ws.Range("F35").Value = 10907200
$ws.Range("G35").Value = 12652500
$ws.Range("H35").Value = 4017800
$ws.Range("I35").Value = 2662700
$ws.Range("J35").Value = 3979700
$ws.Range("D41").Value = 16123800
$ws.Range("E41").Value = 8828900
$ws.Range("F41").Value = 9746400
$ws.Range("G41").Value = 11632200
$ws.Range("H41").Value = 7335300
$ws.Range("I41").Value = 5674900
$ws.Range("J41").Value = 6148000
$ws.Range("D42").Value = 41700
$ws.Range("E42").Value = 82400
$ws.Range("F42").Value = 49000
$ws.Range("G42").Value = 74300
$ws.Range("H42").Value = 137100
$ws.Range("I42").Value = 45900
$ws.Range("J42").Value = 5541800
$ws.Range("D43").Value = 1843600
$ws.Range("E43").Value = 2035900
$ws.Range("F43").Value = 1329400
$ws.Range("G43").Value = 1378100
$ws.Range("H43").Value = 1488600
$ws.Range("I43").Value = 1482100
$ws.Range("J43").Value = 1432800
$ws.Range("D44").Value = 1252100
$ws.Range("E44").Value = 1311000
$ws.Range("F44").Value = 1168000
$ws.Range("G44").Value = 1007900
$ws.Range("H44").Value = 846600
$ws.Range("I44").Value = 809000
$ws.Range("J44").Value = 1808500
$ws.Range("D45").Value = 270100
$ws.Range("E45").Value = 291800
$ws.Range("F45").Value = 354700
$ws.Range("G45").Value = 653000
$ws.Range("H45").Value = 450900
$ws.Range("I45").Value = 442500
$ws.Range("J45").Value = 554400
$ws.Range("D46").Value = 19531400
$ws.Range("E46").Value = 12550100
$ws.Range("F46").Value = 12647600
$ws.Range("G46").Value = 14745500
$ws.Range("H46").Value = 10258500
$ws.Range("I46").Value = 8454300
$ws.Range("J46").Value = 8798100
$ws.Range("D47").Value = 23541800
$ws.Range("E47").Value = 27126500
$ws.Range("F47").Value = 30600400
$ws.Range("G47").Value = 18840800
$ws.Range("H47").Value = 10145400
$ws.Range("I47").Value = 9062100
$ws.Range("J47").Value = 15075600
$ws.Range("D48").Value = 22033700
$ws.Range("E48").Value = 20506500
$ws.Range("F48").Value = 18954200
$ws.Range("G48").Value = 17132100
$ws.Range("H48").Value = 16136100
$ws.Range("I48").Value = 14175300
$ws.Range("J48").Value = 21977200
$ws.Range("D49").Value = 139800
$ws.Range("E49").Value = 142800
$ws.Range("F49").Value = 133500
$ws.Range("G49").Value = 131900
$ws.Range("H49").Value = 125400
$ws.Range("I49").Value = 128100
$ws.Range("J49").Value = 236200
$ws.Range("D52").Value = 44400
$ws.Range("E52").Value = 66500
$ws.Range("F52").Value = 61500
$ws.Range("G52").Value = 63500
$ws.Range("H52").Value = 63400
$ws.Range("I52").Value = 73000
$ws.Range("J52").Value = 62800
$ws.Range("D54").Value = 65291200
$ws.Range("E54").Value = 60392400
$ws.Range("F54").Value = 62397100
$ws.Range("G54").Value = 50913800
$ws.Range("H54").Value = 36728800
$ws.Range("I54").Value = 31892800
$ws.Range("J54").Value = 29199200
$ws.Range("D57").Value = 159800
$ws.Range("E57").Value = 145000
$ws.Range("F57").Value = 158100
$ws.Range("G57").Value = 156000
$ws.Range("H57").Value = 138400
$ws.Range("I57").Value = 565000
$ws.Range("J57").Value = 1682600
$ws.Range("D58").Value = 2840700
$ws.Range("E58").Value = 1382100
$ws.Range("F58").Value = 1217800
$ws.Range("G58").Value = 580100
$ws.Range("H58").Value = 341200
$ws.Range("I58").Value = 282200
$ws.Range("J58").Value = 234400
$ws.Range("D59").Value = 1985600
$ws.Range("E59").Value = 1743800
$ws.Range("F59").Value = 2029300
$ws.Range("G59").Value = 2002800
$ws.Range("H59").Value = 1380200
$ws.Range("I59").Value = 663800
$ws.Range("J59").Value = 821100
$ws.Range("D60").Value = 4986100
$ws.Range("E60").Value = 3270900
$ws.Range("F60").Value = 3405200
$ws.Range("G60").Value = 2738900
$ws.Range("H60").Value = 1859800
$ws.Range("I60").Value = 1511000
$ws.Range("J60").Value = 1579200
$ws.Range("D61").Value = 59300
$ws.Range("E61").Value = 254300
$ws.Range("F61").Value = 67000
$ws.Range("G61").Value = 203200
$ws.Range("H61").Value = 201300
$ws.Range("J61").Value = 400
$ws.Range("D62").Value = 4602600
$ws.Range("E62").Value = 3807800
$ws.Range("F62").Value = 3752400
$ws.Range("G62").Value = 3172300
$ws.Range("H62").Value = 2902600
$ws.Range("I62").Value = 2456300
$ws.Range("J62").Value = 1884500
$ws.Range("D66").Value = 9652600
$ws.Range("E66").Value = 7336800
$ws.Range("F66").Value = 7227500
$ws.Range("G66").Value = 6117100
$ws.Range("H66").Value = 4966400
$ws.Range("I66").Value = 3969900
$ws.Range("J66").Value = 3465800
$ws.Range("D72").Value = 52366900
$ws.Range("E72").Value = 49783900
$ws.Range("F72").Value = 51897800
$ws.Range("G72").Value = 41525000
$ws.Range("H72").Value = 28492000
$ws.Range("I72").Value = 24652600
$ws.Range("J72").Value = 22505800
$ws.Range("D76").Value = 55638600
$ws.Range("E76").Value = 53055600
$ws.Range("F76").Value = 55169600
$ws.Range("G76").Value = 44796700
$ws.Range("H76").Value = 31762400
$ws.Range("I76").Value = 27923000
$ws.Range("J76").Value = 25733400
$ws.Range("D81").Value = 2997700
$ws.Range("E81").Value = -956300
$ws.Range("F81").Value = 10907200
$ws.Range("G81").Value = 12652500
$ws.Range("H81").Value = 4017800
$ws.Range("I81").Value = 2662700
$ws.Range("J81").Value = 3979700
$ws.Range("D83").Value = 1021000
$ws.Range("E83").Value = 1187600
$ws.Range("F83").Value = 1087600
$ws.Range("G83").Value = 1149000
$ws.Range("H83").Value = 746200
$ws.Range("I83").Value = 729100
$ws.Range("J83").Value = 620900
$ws.Range("D89").Value = 5483700
$ws.Range("E89").Value = 4041500
$ws.Range("F89").Value = 3064700
$ws.Range("G89").Value = 2790200
$ws.Range("H89").Value = 3965700
$ws.Range("I89").Value = 3833400
$ws.Range("J89").Value = 3694400
$ws.Range("D91").Value = -2466800
$ws.Range("E91").Value = -2784100
$ws.Range("F91").Value = -2631100
$ws.Range("G91").Value = -2450900
$ws.Range("H91").Value = -2722600
$ws.Range("I91").Value = -2151400
$ws.Range("J91").Value = -1939700
$ws.Range("D94").Value = -4470600
$ws.Range("E94").Value = -3095100
$ws.Range("F94").Value = -956100
$ws.Range("G94").Value = -2594700
$ws.Range("H94").Value = -3742800
$ws.Range("I94").Value = -2980100
$ws.Range("J94").Value = -3296500
$ws.Range("D96").Value = -401800
$ws.Range("E96").Value = -1143000
$ws.Range("F96").Value = -1320400
$ws.Range("G96").Value = -607000
$ws.Range("H96").Value = -448600
$ws.Range("I96").Value = -581600
$ws.Range("J96").Value = -418100
$ws.Range("D100").Value = 677200
$ws.Range("E100").Value = -865400
$ws.Range("F100").Value = -939100
$ws.Range("G100").Value = -436300
$ws.Range("H100").Value = -348600
$ws.Range("I100").Value = -489800
$ws.Range("J100").Value = -487100
$ws.Range("E101").Value = -53900
$ws.Range("F101").Value = 28300
$ws.Range("G101").Value = 162700
$ws.Range("H101").Value = 3000
$ws.Range("I101").Value = -25400
$ws.Range("J101").Value = 6100
$ws.Range("D102").Value = 1694700
$ws.Range("E102").Value = 27100
$ws.Range("F102").Value = 1197800
$ws.Range("G102").Value = -78100
$ws.Range("H102").Value = -122800
$ws.Range("I102").Value = 338100
$ws.Range("J102").Value = -83100
